# Repro: update the Surface accessories link to the general Surface page,
# fill A2/A3 with the same link text and hyperlink them too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl = "https://www.microsoft.com/en-us/surface"

# A1 already carries the old URL text/shared-string; just re-point it.
# A2/A3 need the display text written in first so they pick up the same
# shared string as A1.
$ws.Range("A1").Value = $newUrl
$ws.Range("A2").Value = $newUrl
$ws.Range("A3").Value = $newUrl

# A3 is a brand-new cell - give it the same "Hyperlink" look A1/A2 already have.
$ws.Range("A3").Style = "Hyperlink"

# Hyperlink A1 on its own (display text defaults to the cell's own text).
$ws.Hyperlinks.Add($ws.Range("A1"), $newUrl)

# Hyperlink A2:A3 together, explicit display text matching the shown URL.
$ws.Hyperlinks.Add($ws.Range("A2:A3"), $newUrl, "", "", $newUrl)

# Adding hyperlinks re-stamps a formatting style onto the cells; put the
# plain "Hyperlink" cell style back so A1:A3 stay on the original style.
$ws.Range("A1").Style = "Hyperlink"
$ws.Range("A2").Style = "Hyperlink"
$ws.Range("A3").Style = "Hyperlink"

# Leave the same cell selected as the saved file (A7).
$ws.Range("A7").Select()
